# The edit renames the internal drawing names of the three logo pictures
# that live in the document's headers/footers:
#   - Footer (default), Pearson logo  : image2.png -> image1.png
#   - Footer (first page), Pearson logo: image2.png -> image1.png
#   - Header (first page), BTEC logo   : image1.jpg -> image2.jpg
#
# InlineShape has no writable "Name" in the Word object model, so each
# picture is temporarily converted to a floating Shape (which does expose
# .Name), renamed, and converted back to an inline picture in place.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Default footer: Pearson Edexcel logo (wp:docPr id="1") ---
$ftrDefault = $sec.Footers.Item(1)
$pearsonDefault = $ftrDefault.Range.InlineShapes.Item(1)
$pearsonDefaultShape = $pearsonDefault.ConvertToShape()
$pearsonDefaultShape.Name = "image1.png"
$pearsonDefaultShape.ConvertToInlineShape()

# --- First-page footer: Pearson Edexcel logo (wp:docPr id="2") ---
$ftrFirst = $sec.Footers.Item(2)
$pearsonFirst = $ftrFirst.Range.InlineShapes.Item(1)
$pearsonFirstShape = $pearsonFirst.ConvertToShape()
$pearsonFirstShape.Name = "image1.png"
$pearsonFirstShape.ConvertToInlineShape()

# --- First-page header: BTEC logo (wp:docPr id="3") ---
$hdrFirst = $sec.Headers.Item(2)
$btec = $hdrFirst.Range.InlineShapes.Item(1)
$btecShape = $btec.ConvertToShape()
$btecShape.Name = "image2.jpg"
$btecShape.ConvertToInlineShape()
